$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("G14").Value = "MyRow*"
$ws1.Range("H14").Value = "UserData*"
$ws1.Range("J14").Value = "someListik*"
$ws1.Range("B15").Value = "MyRow*"
$ws1.Range("B16").Value = "UserData*"
$ws1.Range("B18").Value = "someListik*"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B7").Value = "Step1*"
$ws2.Range("G7").Value = "Row1*"
$ws2.Range("B8").Value = "Step2*"
$ws2.Range("G8").Value = "Row2*"
$ws2.Range("G32").Value = "Calculation*"
$ws2.Range("H32").Value = "anycell*"
$ws2.Range("B33").Value = "Step1*"
$ws2.Range("F33").Value = "Step1*"
$ws2.Range("B34").Value = "Step2*"
$ws2.Range("F34").Value = "Step2*"
$ws2.Range("C54").Value = "Calculation*"
$ws2.Range("B55").Value = "Step1*"
$ws2.Range("B56").Value = "Step2*"
$ws2.Range("C72").Value = "Calculation*"
$ws2.Range("D72").Value = "anycell*"
$ws2.Range("B73").Value = "Step1*"

$ws2.Range("D82:H82").BorderAround(1)

$ws2.Range("D85").Value = 2
$ws2.Range("E85").Value = 2
